$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.529.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +12.78%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.837.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +9.73%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.551"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.95%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.74"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.23%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.24"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.286"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0677"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.99%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0933"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.15%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.097.70"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.858.60"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +11.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.653"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.59%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.40%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.481.53"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +12.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.29"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "261.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0757"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.70"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.80"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.80%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.08%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.83%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.04%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +12.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0523"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.53%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.62"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.553.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.60%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.21%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.649"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.15%  "

# Row 39
$ws.Range("B39").Value = "MinaProtocolToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.29"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +212.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0191"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "85.39"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.35%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.919"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +9.57%  "

# Row 44
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.94%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0528"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.62%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.32%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.58"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +28.94%  "

# Row 49
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.989.15"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +9.80%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.95%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.16%  "
